$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend with two new header cells P1, Q1 ---
# Copy formatting from the existing last header cell (O1) onto the new
# header cells so they pick up the same bold/border/center style (s="1").
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
# For every data row: swap the I/K/M/O values and append two new
# unformatted data cells P (=2) and Q (=2).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2 (new column)
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2 (new column)
}
